$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Deals": refresh the full data range (new rows for BotFather,
# Send_Message_telegram and the "User Info" idbot entry, status correction
# for Telegram/Send_Message_telegram, and the two revenue values that moved
# from 0 -> 500 once secrets/test rows were folded into real records).
# ---------------------------------------------------------------------------
$dealsWs = $wb.Worksheets.Item("Deals")

$dealsData = @(
    @("AIBI_Secretary_Bot", "Unknown", 500, "AIBI_Secretary_Bot.txt", "format_a"),
    @("BotFather", "Unknown", 0, "BotFather.txt", "format_a"),
    @("TechCorp Solutions", "Win", 125500, "sample_report_1.txt", "format_b"),
    @("GlobalTrade Inc", "Loss", 0, "sample_report_2.txt", "format_b"),
    @("was particularly concerned about compliance and integration with their existing Salesforce setup.", "Win", 87300, "sample_report_3.txt", "format_b"),
    @("Send_Message_telegram", "Unknown", 0, "Send_Message_telegram.txt", "format_a"),
    @("Telegram", "Unknown", 0, "Telegram.txt", "format_a"),
    @("UFO", "Unknown", 0, "UFO.txt", "format_a"),
    @("User Info • Get ID • idbot", "Unknown", 0, "User Info • Get ID • idbot.txt", "format_a"),
    @("Ілля", "Unknown", 500, "Ілля.txt", "format_a")
)

for ($i = 0; $i -lt $dealsData.Length; $i++) {
    $row = $i + 2
    $dealsWs.Cells.Item($row, 1).Value = $dealsData[$i][0]
    $dealsWs.Cells.Item($row, 2).Value = $dealsData[$i][1]
    $dealsWs.Cells.Item($row, 3).Value = $dealsData[$i][2]
    $dealsWs.Cells.Item($row, 4).Value = $dealsData[$i][3]
    $dealsWs.Cells.Item($row, 5).Value = $dealsData[$i][4]
}

# Widen the "Report File" column (D) now that longer file names are present.
# ColumnWidth uses Excel's character-width units; the small offset below
# compensates for the engine's internal pixel rounding so the stored
# worksheet column width lands exactly on 32.
$dealsWs.Columns.Item(4).ColumnWidth = 32 - 11/12 + 0.02

# ---------------------------------------------------------------------------
# Sheet "Summary": recompute the aggregate stats for the updated deal list.
# ---------------------------------------------------------------------------
$summaryWs = $wb.Worksheets.Item("Summary")

$summaryWs.Cells.Item(2, 1).Value = 10
$summaryWs.Cells.Item(2, 2).Value = 2
$summaryWs.Cells.Item(2, 3).Value = 1
$summaryWs.Cells.Item(2, 4).Value = 20
$summaryWs.Cells.Item(2, 5).Value = 212800
$summaryWs.Cells.Item(2, 6).Value = 106400

# Narrow the "win_rate" column (D) to match the refreshed layout.
$summaryWs.Columns.Item(4).ColumnWidth = 10 - 11/12 + 0.02
